{"js": "// Update the email address shown in the CONTACT section.\nconst body = context.document.body;\nconst results = body.search(\"camilax.gomez97@gmail.com\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"camila.cp.gomez@gmail.com\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the email address shown in the CONTACT section.\n$d = $word.ActiveDocument\n\n$oldEmail = \"camilax.gomez97@gmail.com\"\n$newEmail = \"camila.cp.gomez@gmail.com\"\n\n$target = $null\nforeach ($para in $d.Paragraphs) {\n    if ($para.Range.Text -like \"*$oldEmail*\") {\n        $target = $para\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $r = $target.Range\n    # Exclude the trailing paragraph mark from the replacement range.\n    $r.End = $r.End - 1\n    $r.Text = $newEmail\n}\n"}
